$wb = $excel.ActiveWorkbook

# --- Insert a new worksheet "time_log" between "Sheet1" and "Totals" ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$totals = $wb.Worksheets.Item("Totals")
$newSheet = $wb.Worksheets.Add($totals)
$newSheet.Name = "time_log"

# Headers (order matters for shared-string table interning order)
$newSheet.Range("A1").Value = "date"
$newSheet.Range("C1").Value = "time_spent"
$newSheet.Range("B1").Value = "largest_category"

# Data row
$newSheet.Range("A2").Value = (Get-Date -Year 2023 -Month 5 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$sheet1.Range("C4").Copy()
$newSheet.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$newSheet.Range("B2").Value = "run_me clean, debug, & run"
$newSheet.Range("C2").Value = 4.57

# --- Make time_log the active/selected sheet & cell ---
$newSheet.Activate()
$newSheet.Range("G14").Select()

# --- Update Sheet1 view/formatting ---
$sheet1.Activate()
$sheet1.Range("F8").Select()
$excel.ActiveWindow.Zoom = 100

$sheet1.Columns.Item(3).ColumnWidth = 12.5

# Re-activate time_log as the final active tab
$newSheet.Activate()
